$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44162
$ws.Range("L2").Value = "Tercera"
$ws.Range("M2").Value = 500
$ws.Range("N2").Value = 15000
$ws.Range("O2").Value = 16000
$ws.Range("P2").Value = 15500
$ws.Range("Q2").Value = "$/caja 15 kilos"
$ws.Range("R2").Value = "Región de O'Higgins"
$ws.Range("S2").Value = 1033
$ws.Range("T2").Value = 15

# Row 5
$ws.Range("D5").Value = 44523
$ws.Range("L5").Value = "Segunda"
$ws.Range("M5").Value = 500
$ws.Range("N5").Value = 28000
$ws.Range("O5").Value = 28500
$ws.Range("P5").Value = 28250
$ws.Range("Q5").Value = "$/caja 18 kilos"
$ws.Range("R5").Value = "Provincia de Limarí"
$ws.Range("S5").Value = 1569
$ws.Range("T5").Value = 18

# Row 6
$ws.Range("D6").Value = 44169
$ws.Range("L6").Value = "Segunda"
$ws.Range("M6").Value = 500
$ws.Range("N6").Value = 15000
$ws.Range("O6").Value = 16000
$ws.Range("P6").Value = 15500
$ws.Range("Q6").Value = "$/caja 15 kilos"
$ws.Range("R6").Value = "Región de O'Higgins"
$ws.Range("S6").Value = 1033
$ws.Range("T6").Value = 15

# Row 7
$ws.Range("D7").Value = 44194
$ws.Range("L7").Value = "Segunda"
$ws.Range("M7").Value = 300
$ws.Range("N7").Value = 15000
$ws.Range("O7").Value = 16000
$ws.Range("P7").Value = 15500
$ws.Range("Q7").Value = "$/caja 15 kilos"
$ws.Range("R7").Value = "Región Metropolitana"
$ws.Range("S7").Value = 1033
$ws.Range("T7").Value = 15

# Row 8
$ws.Range("D8").Value = 44176
$ws.Range("L8").Value = "Segunda"
$ws.Range("M8").Value = 500
$ws.Range("N8").Value = 15000
$ws.Range("O8").Value = 16000
$ws.Range("P8").Value = 15500
$ws.Range("Q8").Value = "$/caja 15 kilos"
$ws.Range("R8").Value = "Región Metropolitana"
$ws.Range("S8").Value = 1033
$ws.Range("T8").Value = 15

# Row 9
$ws.Range("D9").Value = 44530
$ws.Range("L9").Value = "Segunda"
$ws.Range("M9").Value = 500
$ws.Range("N9").Value = 20000
$ws.Range("O9").Value = 21000
$ws.Range("P9").Value = 20500
$ws.Range("Q9").Value = "$/caja 18 kilos"
$ws.Range("R9").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S9").Value = 1139
$ws.Range("T9").Value = 18

# Row 10
$ws.Range("D10").Value = 44159
$ws.Range("L10").Value = "Tercera"
$ws.Range("M10").Value = 400
$ws.Range("N10").Value = 15500
$ws.Range("O10").Value = 16000
$ws.Range("P10").Value = 15750
$ws.Range("Q10").Value = "$/caja 15 kilos"
$ws.Range("R10").Value = "Región de O'Higgins"
$ws.Range("S10").Value = 1050
$ws.Range("T10").Value = 15

# Row 11
$ws.Range("D11").Value = 44166
$ws.Range("L11").Value = "Segunda"
$ws.Range("M11").Value = 600
$ws.Range("N11").Value = 16000
$ws.Range("O11").Value = 17000
$ws.Range("P11").Value = 16500
$ws.Range("Q11").Value = "$/caja 15 kilos"
$ws.Range("R11").Value = "Región de O'Higgins"
$ws.Range("S11").Value = 1100
$ws.Range("T11").Value = 15
